# "Alterações nas prioridades das sprints"
#
# Adds a third "responsible" sub-column (E/F) to the Sprint 2 block
# (rows 9-15), mirroring the structure already used by the Sprint 1
# block (columns A-B / C-D) and by the Sprint 3 block (which already
# has a third E/F sub-column). Also updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# --- Copy formatting from matching "template" cells so the new
#     E/F columns look identical to their siblings -------------------

# Header row (row 9): E9 should look like the other wrap-text category
# headers (e.g. E17 "Funcionários Cadastrados"), F9 like the other
# "Responsável" headers (e.g. F17).
$ws.Range("E17").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F17").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null

# Data rows 10-12: mirror the formatting already used in columns C/D
# of the same row.
$ws.Range("C10").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").Copy() | Out-Null
$ws.Range("F10").PasteSpecial(-4122) | Out-Null

$ws.Range("C11").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null

$ws.Range("C12").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null

# Rows 13-15: only column E picks up the formatting of the (empty)
# column C cells in those rows - column F is left as-is.
$ws.Range("C13").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Fill in the new values ------------------------------------------

$ws.Range("E9").Value  = "Clientes Cadastrados"
$ws.Range("F9").Value  = "Responsável"

$ws.Range("E10").Value = "Mockaps"
$ws.Range("F10").Value = "vinicuios"

$ws.Range("E11").Value = "Desenvolvida"
$ws.Range("F11").Value = "Guilherme"

$ws.Range("E12").Value = "Homologada"
$ws.Range("F12").Value = "Natan"

# --- Update the active cell / selection -------------------------------

$ws.Range("F13").Select() | Out-Null
